$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: header text updates (shared strings) ---
$ws.Range("A8").Value = "Volume 31   Number  25"
$ws.Range("C9").Value = "Report Covering the Week  6/17/2024  Through  6/23/2024"

# --- Step 2: cells that change type (number <-> blank-marker text) ---
# Blank-marker cells reuse the shared strings "0" / "***.*" with the
# General-format style (style index 14). Copying a donor cell that
# already has the desired style+value reproduces both the style index
# and the shared-string reference exactly, matching native Excel output.
$ws.Range("C14").Copy($ws.Range("G15"))
$ws.Range("E14").Copy($ws.Range("H15"))
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("I14").Copy($ws.Range("C28"))
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("C14").Copy($ws.Range("D31"))
$ws.Range("E14").Copy($ws.Range("E31"))
$ws.Range("I14").Copy($ws.Range("F31"))

# --- Step 3: plain numeric value updates ---
$ws.Range("L15").Value = -8.333333333333
$ws.Range("M15").Value = 120
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 10
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 60
$ws.Range("J16").Value = 63
$ws.Range("K16").Value = -4.761904761904
$ws.Range("L16").Value = -26.829268292682
$ws.Range("M16").Value = -53.125
$ws.Range("N16").Value = -86.547085201793
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 100
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = -5.263157894736
$ws.Range("I17").Value = 135
$ws.Range("J17").Value = 104
$ws.Range("K17").Value = 29.807692307692
$ws.Range("L17").Value = 58.823529411764
$ws.Range("M17").Value = 229.268292682927
$ws.Range("N17").Value = 27.358490566037
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 33.333333333333
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 107
$ws.Range("J18").Value = 130
$ws.Range("K18").Value = -17.692307692307
$ws.Range("L18").Value = 10.309278350515
$ws.Range("M18").Value = -10.084033613445
$ws.Range("N18").Value = -85.241379310344
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = -8.108108108108
$ws.Range("I19").Value = 267
$ws.Range("J19").Value = 282
$ws.Range("K19").Value = -5.319148936170
$ws.Range("L19").Value = -14.696485623003
$ws.Range("M19").Value = 13.617021276595
$ws.Range("N19").Value = -7.931034482758
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 29
$ws.Range("G20").Value = 50
$ws.Range("H20").Value = -42
$ws.Range("I20").Value = 155
$ws.Range("J20").Value = 158
$ws.Range("K20").Value = -1.898734177215
$ws.Range("L20").Value = 32.478632478632
$ws.Range("M20").Value = 58.163265306122
$ws.Range("N20").Value = -93.341924398625
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -7.407407407407
$ws.Range("F21").Value = 106
$ws.Range("G21").Value = 130
$ws.Range("H21").Value = -18.461538461538
$ws.Range("I21").Value = 738
$ws.Range("J21").Value = 745
$ws.Range("K21").Value = -0.939597315436
$ws.Range("L21").Value = 4.532577903682
$ws.Range("M21").Value = 17.515923566879
$ws.Range("N21").Value = -81.144609095554
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -80
$ws.Range("M22").Value = 38.461538461538
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -66.666666666666
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = -42.857142857142
$ws.Range("I23").Value = 24
$ws.Range("J23").Value = 37
$ws.Range("K23").Value = -35.135135135135
$ws.Range("L23").Value = -11.111111111111
$ws.Range("M23").Value = 33.333333333333
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = -36
$ws.Range("F24").Value = 87
$ws.Range("G24").Value = 105
$ws.Range("H24").Value = -17.142857142857
$ws.Range("I24").Value = 572
$ws.Range("J24").Value = 707
$ws.Range("K24").Value = -19.094766619519
$ws.Range("L24").Value = -10.625
$ws.Range("M24").Value = 25.164113785558
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -22.222222222222
$ws.Range("F25").Value = 32
$ws.Range("G25").Value = 40
$ws.Range("H25").Value = -20
$ws.Range("I25").Value = 257
$ws.Range("J25").Value = 261
$ws.Range("K25").Value = -1.532567049808
$ws.Range("L25").Value = 39.673913043478
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 40
$ws.Range("F26").Value = 44
$ws.Range("G26").Value = 37
$ws.Range("H26").Value = 18.918918918918
$ws.Range("I26").Value = 242
$ws.Range("J26").Value = 227
$ws.Range("K26").Value = 6.607929515418
$ws.Range("L26").Value = 10.502283105022
$ws.Range("M26").Value = 19.801980198019
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 0
$ws.Range("L27").Value = -26.315789473684
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = 33.333333333333
$ws.Range("I28").Value = 21
$ws.Range("K28").Value = -4.545454545454
$ws.Range("L28").Value = -34.375
$ws.Range("G31").Value = 3
$ws.Range("H31").Value = -66.666666666666
$ws.Range("I31").Value = 2
$ws.Range("K31").Value = -86.666666666666
$ws.Range("L31").Value = 0

